$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-09-14 Thursday" "2023-09-15 Friday"

Replace-Text "90×93=" "84×39="
Replace-Text "17×41=" "42×25="
Replace-Text "61×34=" "14×48="
Replace-Text "99×43=" "20×11="
Replace-Text "91×68=" "17×32="

Replace-Text "39×57=" "82×91="
Replace-Text "89×31=" "48×52="
Replace-Text "68×73=" "86×78="
Replace-Text "29×77=" "60×76="
Replace-Text "92×54=" "26×66="

Replace-Text "18×33=" "59×51="
Replace-Text "16×95=" "60×17="
Replace-Text "44×47=" "71×95="
Replace-Text "72×40=" "88×33="
Replace-Text "87×25=" "46×95="

Replace-Text "21×27=" "46×67="
Replace-Text "46×12=" "88×72="
Replace-Text "69×69=" "85×37="
Replace-Text "88×79=" "32×60="
Replace-Text "30×80=" "28×84="

Replace-Text "97×66=" "77×46="
Replace-Text "60×88=" "17×35="
Replace-Text "98×95=" "94×92="
Replace-Text "37×66=" "44×76="
Replace-Text "57×39=" "54×51="
